$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure date-like text cells (Y, AA columns) stay as text, not auto-converted to dates
$ws.Range("Y2:Y24").NumberFormat = "@"
$ws.Range("AA2:AA24").NumberFormat = "@"

$ws.Range("A2").Value = 111784168
$ws.Range("B2").Value = 90709
$ws.Range("E2").Value = 5448
$ws.Range("F2").Value = "Svartvit taggsvamp"
$ws.Range("G2").Value = "Phellodon connatus"
$ws.Range("H2").Value = "(Schultz) nom.prov"
$ws.Range("Q2").Value = 504976.0557203053
$ws.Range("R2").Value = 7018779.864305317
$ws.Range("Y2").Value = "2023-08-23"
$ws.Range("AA2").Value = "2023-08-23"
$ws.Range("A3").Value = 111783769
$ws.Range("B3").Value = 90665
$ws.Range("D3").Value = "VU"
$ws.Range("E3").Value = 1435
$ws.Range("F3").Value = "Bitter taggsvamp"
$ws.Range("G3").Value = "Hydnellum fennicum"
$ws.Range("H3").Value = "(P.Karst.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q3").Value = 505036.7939151306
$ws.Range("R3").Value = 7018819.987804689
$ws.Range("A4").Value = 111782876
$ws.Range("B4").Value = 90651
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 1968
$ws.Range("F4").Value = "Grantaggsvamp"
$ws.Range("G4").Value = "Bankera violascens"
$ws.Range("H4").Value = "(Alb. & Schwein. : Fr.) Pouzar"
$ws.Range("Q4").Value = 505022.9813479512
$ws.Range("R4").Value = 7018724.615566149
$ws.Range("A5").Value = 111783071
$ws.Range("B5").Value = 90300
$ws.Range("E5").Value = 4745
$ws.Range("F5").Value = "Tallriska"
$ws.Range("G5").Value = "Lactarius musteus"
$ws.Range("H5").Value = "Fr."
$ws.Range("Q5").Value = 505060.2648977584
$ws.Range("R5").Value = 7018787.191973396
$ws.Range("A6").Value = 111782750
$ws.Range("B6").Value = 82949
$ws.Range("E6").Value = 5589
$ws.Range("F6").Value = "Rödbrun klubbdyna"
$ws.Range("G6").Value = "Trichoderma nybergianum"
$ws.Range("H6").Value = "(T.Ulvinen & H.L.Chamb.) Jaklitsch & Voglmayr"
$ws.Range("Q6").Value = 505007.618534557
$ws.Range("R6").Value = 7018756.52538473
$ws.Range("Y6").Value = "2023-08-30"
$ws.Range("AA6").Value = "2023-08-30"
$ws.Range("A7").Value = 111814348
$ws.Range("B7").Value = 96370
$ws.Range("E7").Value = 219847
$ws.Range("F7").Value = "Tvåblad"
$ws.Range("G7").Value = "Neottia ovata"
$ws.Range("H7").Value = "(L.) Buff. & Fingerh."
$ws.Range("Q7").Value = 504944.9568800884
$ws.Range("R7").Value = 7018794.658574538
$ws.Range("A8").Value = 111814402
$ws.Range("B8").Value = 90651
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 1968
$ws.Range("F8").Value = "Grantaggsvamp"
$ws.Range("G8").Value = "Bankera violascens"
$ws.Range("H8").Value = "(Alb. & Schwein. : Fr.) Pouzar"
$ws.Range("Q8").Value = 505200.3682009591
$ws.Range("R8").Value = 7018764.927175661
$ws.Range("A9").Value = 111814411
$ws.Range("B9").Value = 96253
$ws.Range("E9").Value = 504
$ws.Range("F9").Value = "Guckusko"
$ws.Range("G9").Value = "Cypripedium calceolus"
$ws.Range("H9").Value = "L."
$ws.Range("Q9").Value = 505014.8575873387
$ws.Range("R9").Value = 7018735.397438973
$ws.Range("A10").Value = 111814362
$ws.Range("B10").Value = 103288
$ws.Range("E10").Value = 221144
$ws.Range("F10").Value = "Grönpyrola"
$ws.Range("G10").Value = "Pyrola chlorantha"
$ws.Range("H10").Value = "Sw."
$ws.Range("Q10").Value = 504958.3523041067
$ws.Range("R10").Value = 7018869.788911887
$ws.Range("A11").Value = 111814344
$ws.Range("B11").Value = 103288
$ws.Range("D11").Value = "LC"
$ws.Range("E11").Value = 221144
$ws.Range("F11").Value = "Grönpyrola"
$ws.Range("G11").Value = "Pyrola chlorantha"
$ws.Range("H11").Value = "Sw."
$ws.Range("Q11").Value = 505070.3462143789
$ws.Range("R11").Value = 7018677.46330901
$ws.Range("A12").Value = 111814415
$ws.Range("B12").Value = 56543
$ws.Range("E12").Value = 103021
$ws.Range("F12").Value = "Talltita"
$ws.Range("G12").Value = "Poecile montanus"
$ws.Range("H12").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("Q12").Value = 505100.2035408606
$ws.Range("R12").Value = 7018878.55609256
$ws.Range("AJ12").ClearContents()
$ws.Range("AK12").ClearContents()
$ws.Range("AO12").ClearContents()
$ws.Range("A13").Value = 111814428
$ws.Range("B13").Value = 90666
$ws.Range("D13").Value = "LC"
$ws.Range("E13").Value = 4364
$ws.Range("F13").Value = "Dropptaggsvamp"
$ws.Range("G13").Value = "Hydnellum ferrugineum"
$ws.Range("H13").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q13").Value = 504987.8220338543
$ws.Range("R13").Value = 7018743.451279385
$ws.Range("A14").Value = 111814417
$ws.Range("Q14").Value = 505110.448201828
$ws.Range("R14").Value = 7018666.73204405
$ws.Range("A15").Value = 111814369
$ws.Range("B15").Value = 96265
$ws.Range("D15").Value = "LC"
$ws.Range("E15").Value = 219790
$ws.Range("F15").Value = "Fläcknycklar"
$ws.Range("G15").Value = "Dactylorhiza maculata"
$ws.Range("H15").Value = "(L.) Soó"
$ws.Range("Q15").Value = 504944.0561301867
$ws.Range("R15").Value = 7018794.657189432
$ws.Range("A16").Value = 111814356
$ws.Range("B16").Value = 90678
$ws.Range("E16").Value = 4366
$ws.Range("F16").Value = "Skarp dropptaggsvamp"
$ws.Range("G16").Value = "Hydnellum peckii"
$ws.Range("H16").Value = "Banker"
$ws.Range("Q16").Value = 505204.4099656619
$ws.Range("R16").Value = 7018772.129998797
$ws.Range("A17").Value = 111814388
$ws.Range("B17").Value = 78578
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 6458
$ws.Range("F17").Value = "Lunglav"
$ws.Range("G17").Value = "Lobaria pulmonaria"
$ws.Range("H17").Value = "(L.) Hoffm."
$ws.Range("Q17").Value = 505182.7410700406
$ws.Range("R17").Value = 7018803.578552675
$ws.Range("AJ17").Value = "sälg"
$ws.Range("AK17").Value = "Salix caprea"
$ws.Range("AO17").Value = "Salix caprea"
$ws.Range("A18").Value = 111814375
$ws.Range("B18").Value = 96370
$ws.Range("E18").Value = 219847
$ws.Range("F18").Value = "Tvåblad"
$ws.Range("G18").Value = "Neottia ovata"
$ws.Range("H18").Value = "(L.) Buff. & Fingerh."
$ws.Range("Q18").Value = 504999.9977373667
$ws.Range("R18").Value = 7018733.575208749
$ws.Range("A19").Value = 111814434
$ws.Range("B19").Value = 96253
$ws.Range("D19").Value = "LC"
$ws.Range("E19").Value = 504
$ws.Range("F19").Value = "Guckusko"
$ws.Range("G19").Value = "Cypripedium calceolus"
$ws.Range("H19").Value = "L."
$ws.Range("Q19").Value = 505012.6056710624
$ws.Range("R19").Value = 7018735.393927739
$ws.Range("A20").Value = 111814395
$ws.Range("B20").Value = 103288
$ws.Range("D20").Value = "LC"
$ws.Range("E20").Value = 221144
$ws.Range("F20").Value = "Grönpyrola"
$ws.Range("G20").Value = "Pyrola chlorantha"
$ws.Range("H20").Value = "Sw."
$ws.Range("Q20").Value = 504987.8220338543
$ws.Range("R20").Value = 7018743.451279385
$ws.Range("A21").Value = 111814432
$ws.Range("B21").Value = 95674
$ws.Range("E21").Value = 222741
$ws.Range("F21").Value = "Finbräken"
$ws.Range("G21").Value = "Cystopteris montana"
$ws.Range("H21").Value = "(Lam.) Desv."
$ws.Range("Q21").Value = 505015.75484597
$ws.Range("R21").Value = 7018737.647699019
$ws.Range("A22").Value = 111814359
$ws.Range("B22").Value = 90300
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 4745
$ws.Range("F22").Value = "Tallriska"
$ws.Range("G22").Value = "Lactarius musteus"
$ws.Range("H22").Value = "Fr."
$ws.Range("Q22").Value = 505073.4975346876
$ws.Range("R22").Value = 7018678.36784017
$ws.Range("A23").Value = 111814350
$ws.Range("B23").Value = 89401
$ws.Range("D23").Value = "NT"
$ws.Range("E23").Value = 1108
$ws.Range("F23").Value = "Harticka"
$ws.Range("G23").Value = "Pelloporus leporinus"
$ws.Range("H23").Value = "(Fr.) Krieglst."
$ws.Range("Q23").Value = 504984.8875472886
$ws.Range("R23").Value = 7018893.217038274
$ws.Range("A24").Value = 111814351
$ws.Range("B24").Value = 89845
$ws.Range("D24").Value = "VU"
$ws.Range("E24").Value = 1209
$ws.Range("F24").Value = "Rynkskinn"
$ws.Range("G24").Value = "Phlebia centrifuga"
$ws.Range("H24").Value = "P.Karst."
$ws.Range("Q24").Value = 504984.8875472886
$ws.Range("R24").Value = 7018893.217038274
